# Commit: updated RegisterPage first name By locator, updated product sheet in Excel sheet
#
# 1) Rename the existing sheet to "register"
# 2) Store the telephone numbers as text (quote-prefixed) instead of numbers
# 3) Refresh the mailto hyperlinks so they no longer carry an explicit display string
# 4) Add a new "product" worksheet after "register" with search-key / product-name data
# 5) Restore view state (active sheet, selection, zoom) on both sheets
# 6) Tighten up column widths to better fit the new content

$wb = $excel.ActiveWorkbook

# --- register sheet -------------------------------------------------------
$register = $wb.Worksheets.Item(1)
$register.Name = "register"

# Telephone numbers become text values (quotePrefix) rather than numbers
$register.Range("C2").Value = "'7878789898"
$register.Range("C3").Value = "'7878789890"
$register.Range("C4").Value = "'7878789843"

# Recreate the hyperlinks without an explicit display override
$register.Hyperlinks.Delete()
$register.Hyperlinks.Add($register.Range("D2"), "mailto:harpreet@123")
$register.Hyperlinks.Add($register.Range("D3"), "mailto:ratul@123")
$register.Hyperlinks.Add($register.Range("D4"), "mailto:sandhya@123")

# Column widths tuned for the refreshed layout
$register.Columns.Item(1).ColumnWidth = 10.166666666666666
$register.Columns.Item(2).ColumnWidth = 11.166666666666666
$register.Columns.Item(3).ColumnWidth = 12.330729166666666
$register.Columns.Item(4).ColumnWidth = 11.998697916666666

# --- product sheet (new) ---------------------------------------------------
$product = $wb.Worksheets.Add($null, $register)
$product.Name = "product"

$product.Cells.Item(1, 1).Value = "searchkey"
$product.Cells.Item(1, 2).Value = "productname"
$product.Cells.Item(2, 1).Value = "macbook"
$product.Cells.Item(2, 2).Value = "MacBook Pro"
$product.Cells.Item(3, 1).Value = "samsung"
$product.Cells.Item(4, 1).Value = "imac"
$product.Cells.Item(5, 1).Value = "canon"
$product.Cells.Item(3, 2).Value = "Samsung SyncMaster 941BW"
$product.Cells.Item(4, 2).Value = "iMac"
$product.Cells.Item(5, 2).Value = "Canon EOS 5D"

$product.Columns.Item(1).ColumnWidth = 8.5
$product.Columns.Item(2).ColumnWidth = 23.75

# view state for the product sheet
$product.Range("B16").Select()
$excel.ActiveWindow.Zoom = 173

# --- restore view state on register sheet ----------------------------------
$register.Activate()
$register.Range("E4").Select()
$excel.ActiveWindow.Zoom = 187
